$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "Through 2021-11-26"

# Update the label in A13
$ws.Range("A13").Value = "November (through 11-26)"

# Row 13 updates
$ws.Range("C13").Value = 27
$ws.Range("D13").Value = 0.0357
$ws.Range("F13").Value = 59
$ws.Range("G13").Value = 0.1061
$ws.Range("I13").Value = 96
$ws.Range("J13").Value = 0.0204
$ws.Range("K13").Value = 8
$ws.Range("L13").Value = 50
$ws.Range("M13").Value = 0.1379
$ws.Range("O13").Value = 40
$ws.Range("P13").Value = 0.1304
$ws.Range("R13").Value = 174
$ws.Range("S13").Value = 0.044
$ws.Range("U13").Value = 173
$ws.Range("V13").Value = 0.0226

# Row 14 updates
$ws.Range("C14").Value = 253
$ws.Range("D14").Value = 0.1154
$ws.Range("F14").Value = 493
$ws.Range("G14").Value = 0.1069
$ws.Range("I14").Value = 745
$ws.Range("J14").Value = 0.078
$ws.Range("K14").Value = 74
$ws.Range("L14").Value = 599
$ws.Range("M14").Value = 0.11
$ws.Range("O14").Value = 474
$ws.Range("P14").Value = 0.1023
$ws.Range("R14").Value = 1178
$ws.Range("S14").Value = 0.0492
$ws.Range("U14").Value = 1524
$ws.Range("V14").Value = 0.0593
